$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the newly logged work rows (31-44) ---
# Text cells are written in the same order the shared-string table records
# them so the underlying sharedStrings.xml matches the original edit.

$ws.Range("C31").Value = "studied chapter 19"
$ws.Range("D32").Value = "another doctors appointment"
$ws.Range("C33").Value = "additional research via YT and w3schools"
$ws.Range("C34").Value = "additional research via YT and w3schools"
$ws.Range("C37").Value = "started working on chapter 20"
$ws.Range("C41").Value = "revising "
$ws.Range("D41").Value = "looked up the game project in the book to get some concepts for the practical part"
$ws.Range("D42").Value = "freeday, was moving back to Oulu"
$ws.Range("D35").Value = "freeday was moving from Oulu"
$ws.Range("C44").Value = "documentation work, reorganising working schedule, revision"
$ws.Range("C40").Value = "additional studying of parcers based on prev exercises"
$ws.Range("C38").Value = "almost finished chapter 20"
$ws.Range("D37").Value = "node programming goes very hard"
$ws.Range("D38").Value = "still  very confused with node coding, ned more info"
$ws.Range("D44").Value = "node still does not work well"

$ws.Range("D36").Value = "freeday"
$ws.Range("D39").Value = "freeday"
$ws.Range("D43").Value = "freeday"

# --- Hours column (plain numbers, no shared strings involved) ---

$ws.Range("B31").Value = 3
$ws.Range("B32").Value = 0
$ws.Range("B33").Value = 2
$ws.Range("B34").Value = 3
$ws.Range("B35").Value = 0
$ws.Range("B36").Value = 0
$ws.Range("B37").Value = 2
$ws.Range("B38").Value = 2
$ws.Range("B39").Value = 0
$ws.Range("B40").Value = 3
$ws.Range("B41").Value = 2
$ws.Range("B42").Value = 0
$ws.Range("B43").Value = 0
$ws.Range("B44").Value = 3

# --- Cosmetic changes: column width and current view/selection ---

$ws.Columns.Item(3).ColumnWidth = 49.5

$ws.Application.ActiveWindow.ScrollRow = 23
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D44").Select()
